# Updated to use CLI for deploying model.
# Insert two explanatory rows into the Variables sheet and tweak the
# BASE_NAME sample value to a placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1 for the "create a variable group" instructions.
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "To create a variable group in Azure DevOps, click Pipelines > Library > + Variable group"

# Insert a new row above the (now shifted) "Variables" header row (row 6)
# for the "add the following variables" instructions.
$ws.Rows.Item(6).Insert()
$ws.Range("A6").Value = "Add the following variables into devopsforai-aml-vg"

# Update the BASE_NAME sample value to a generic placeholder.
$ws.Range("B12").Value = "<specify a unique name>"
